$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-EFormat($row, $donor) {
    $ws.Range($donor).Copy()
    $ws.Range("E$row").PasteSpecial(-4122)
    $ws.Range("E$row").NumberFormat = "0.0"
    $ws.Range("E$row").HorizontalAlignment = -4152
    $ws.Range("E$row").VerticalAlignment = -4108
    $ws.Range("E$row").WrapText = $false
}

# Row 3: header year 2020 (copy D3 formatting, which already matches)
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = 2020

# Row 4: value 22.1, border continuing from header (copy D4 formatting)
Set-EFormat 4 "D4"
$ws.Range("E4").Value = 22.1

# Row 5: blank, plain no-border style
Set-EFormat 5 "D5"

# Row 6: value 52.7
Set-EFormat 6 "D5"
$ws.Range("E6").Value = 52.7

# Row 7: value 5
Set-EFormat 7 "D5"
$ws.Range("E7").Value = 5

# Row 8: blank
Set-EFormat 8 "D5"

# Row 9: value 4.8
Set-EFormat 9 "D5"
$ws.Range("E9").Value = 4.8

# Row 10: value 15.8
Set-EFormat 10 "D5"
$ws.Range("E10").Value = 15.8

# Row 11: value 13.5
Set-EFormat 11 "D5"
$ws.Range("E11").Value = 13.5

# Row 12: value 9.6
Set-EFormat 12 "D5"
$ws.Range("E12").Value = 9.6

# Row 13: value 2.7
Set-EFormat 13 "D5"
$ws.Range("E13").Value = 2.7

# Row 14: value 14.7
Set-EFormat 14 "D5"
$ws.Range("E14").Value = 14.7

# Row 15: value 18.2
Set-EFormat 15 "D5"
$ws.Range("E15").Value = 18.2

# Row 16: value 74
Set-EFormat 16 "D5"
$ws.Range("E16").Value = 74

# Row 17: value 35.1
Set-EFormat 17 "D5"
$ws.Range("E17").Value = 35.1

# Row 18: blank
Set-EFormat 18 "D5"

# Rows 19-23: text "-" (data not available), style based on D26's (fill+border applied) xf
Set-EFormat 19 "D26"
$ws.Range("E19").Value = "-"

Set-EFormat 20 "D26"
$ws.Range("E20").Value = "-"

Set-EFormat 21 "D26"
$ws.Range("E21").Value = "-"

Set-EFormat 22 "D26"
$ws.Range("E22").Value = "-"

Set-EFormat 23 "D26"
$ws.Range("E23").Value = "-"

# Row 24: blank, header-ish row with fill-applied style (no border)
Set-EFormat 24 "A25"

# Rows 25-28: text "-"
Set-EFormat 25 "D26"
$ws.Range("E25").Value = "-"

Set-EFormat 26 "D26"
$ws.Range("E26").Value = "-"

Set-EFormat 27 "D26"
$ws.Range("E27").Value = "-"

Set-EFormat 28 "D26"
$ws.Range("E28").Value = "-"

# Row 29: text "-", bottom-border style (thick bottom row)
Set-EFormat 29 "D29"
$ws.Range("E29").Value = "-"

# Update the active selection to match the recorded view state
$ws.Range("J24").Select()
